# Speakers_template.xlsx edit
# 1. Show user full name rather than login name for all kinds of data.
#    -> The "salesPerson" (login name) placeholder becomes "salesPersonFullName".
# 2. Validate region/department required input (application-side change, not
#    represented in this worksheet beyond re-saving it).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the sales person login-name placeholder with the full-name placeholder.
$ws.Range("D2").Value = '${record.salesPersonFullName}'

# Update the selected/active cell to E2 (matches saved sheet view state).
$ws.Range("E2").Select()
